# convert labour inputs to excel
# Adds a "Labour" section to the Price sheet (manager / permanent / casual
# staff cost inputs) together with matching workbook-level defined names
# and explanatory cell comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Price")

# --- New "Labour" section (rows 61-76, values in column B) ---------------
$ws.Range("A61").Value = "Labour"
$ws.Range("A61").Font.Bold = $true

$ws.Range("A62").Value = "manager_cost"
$ws.Range("B62").Value = 80000

$ws.Range("A64").Value = "permanent_cost"
$ws.Range("B64").Value = 80000

$ws.Range("A66").Value = "permanent_super"
$ws.Range("B66").Value = 0.09

$ws.Range("A68").Value = "permanent_workers_comp"
$ws.Range("B68").Value = 0.035

$ws.Range("A70").Value = "permanent_ls_leave"
$ws.Range("B70").Value = 0.023

$ws.Range("A72").Value = "casual_cost"
$ws.Range("B72").Value = 28

$ws.Range("A74").Value = "casual_super"
$ws.Range("B74").Value = 0.09

$ws.Range("A76").Value = "casual_workers_comp"
$ws.Range("B76").Value = 0.035

# --- Cell comments explaining the new inputs ------------------------------
$ws.Range("A62").AddComment("Michael Young:`n`$/yr") | Out-Null
$ws.Range("A64").AddComment("Michael Young:`n`$/yr`nbefore super") | Out-Null
$ws.Range("A70").AddComment("Michael Young:`nLS leave") | Out-Null

# --- Workbook-level defined names -----------------------------------------
$wb.Names.Add('manager_cost', '=Price!$B$62')
$wb.Names.Add('permanent_cost', '=Price!$B$64')
$wb.Names.Add('permanent_super', '=Price!$B$66')
$wb.Names.Add('permanent_workers_comp', '=Price!$B$68')
$wb.Names.Add('permanent_ls_leave', '=Price!$B$70')
$wb.Names.Add('casual_cost', '=Price!$B$72')
$wb.Names.Add('casual_super', '=Price!$B$74')
$wb.Names.Add('casual_workers_comp', '=Price!$B$76')

# --- Keep the view pointed at the new section -----------------------------
$ws.Select() | Out-Null
$ws.Range("B72").Select() | Out-Null
